$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 48, 49, 50, 51, 53 (row 52 unchanged)
# This represents a cyclic rotation of row contents across columns
# A, B, D, E, F, G, H, P, Q, R (other columns remain with their own row's data).

$targets = @{
    48 = @{ A=111974029; B=88032;  D='VU'; E=6276; F='Goliatmusseron';    G='Tricholoma matsutake';   H='(S.Ito & S.Imai) Singer';               P='Aloppmoarna, Jmt';     Q=439334.7866423383; R=6952296.802153576 }
    49 = @{ A=111974134; B=90658;  D='NT'; E=4361; F='Orange taggsvamp';  G='Hydnellum aurantiacum';  H='(Batsch:Fr.) P.Karst.';                 P='Aloppmoarna i S, Jmt'; Q=439399.8222122483; R=6952207.441512506 }
    50 = @{ A=111974126; B=88032;  D='VU'; E=6276; F='Goliatmusseron';    G='Tricholoma matsutake';   H='(S.Ito & S.Imai) Singer';               P='Aloppmoarna i S, Jmt'; Q=439289.9461055733; R=6952209.002200785 }
    51 = @{ A=111974124; B=90666;  D='LC'; E=4364; F='Dropptaggsvamp';    G='Hydnellum ferrugineum';  H='(Fr.:Fr.) P. Karst.';                   P='Aloppmoarna i S, Jmt'; Q=439276.3867801811; R=6952196.853249942 }
    53 = @{ A=111974133; B=90682;  D='NT'; E=2059; F='Skrovlig taggsvamp';G='Hydnellum scabrosum';    H='(Fr.) E.Larss., K.H.Larss. & Kõljalg';  P='Aloppmoarna i S, Jmt'; Q=439389.9449806474; R=6952220.480550999 }
}

foreach ($row in $targets.Keys) {
    $t = $targets[$row]
    $ws.Range("A$row").Value = $t.A
    $ws.Range("B$row").Value = $t.B
    $ws.Range("D$row").Value = $t.D
    $ws.Range("E$row").Value = $t.E
    $ws.Range("F$row").Value = $t.F
    $ws.Range("G$row").Value = $t.G
    $ws.Range("H$row").Value = $t.H
    $ws.Range("P$row").Value = $t.P
    $ws.Range("Q$row").Value = $t.Q
    $ws.Range("R$row").Value = $t.R
}
